# "add file reading practice"
# Fix a couple of typos/labels on the Asia sheet, and replace the
# Africa sheet's sample rows with a fresh set of country/capital pairs.

$wb = $excel.ActiveWorkbook

# --- Asia sheet: Delhi -> New Delhi, South Korea -> South Koera ---
$asia = $wb.Worksheets.Item("Asia")
$asia.Range("B4").Value = "New Delhi"
$asia.Range("C5").Value = "South Koera"

# --- Africa sheet: new city/country sample data ---
$africa = $wb.Worksheets.Item("Africa")
$africa.Range("B2").Value = "Abuja"
$africa.Range("C2").Value = "Nigeria"
$africa.Range("B3").Value = "Cairo"
$africa.Range("C3").Value = "Egypt"
$africa.Range("B4").Value = "Conakry"
$africa.Range("C4").Value = "Guinea"
$africa.Range("B5").Value = "Addis Ababa"
$africa.Range("C5").Value = "Ethiopia"
$africa.Range("D9").Select()

# --- Europe sheet: unchanged data, just note the last selection ---
$europe = $wb.Worksheets.Item("Europe")
$europe.Range("B2").Select()

# --- Asia sheet stays the active tab/selection ---
$asia.Activate()
$asia.Range("B3").Select()
